# DU FBS Mock 7 - enter this student's (row 4, "Sheet1") raw marks.
# The percentage / total / rank columns are formulas and recalculate
# automatically once the raw inputs below are written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 4
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 6
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = 5
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("S4").Value = 13
$ws.Range("T4").Value = 2

# Leave the selection where the author left it after the data entry.
$ws.Range("C4:T4").Select()
